$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = "2.393"
    "D2" = "18"
    "F2" = "188"
    "H2" = "0.033"
    "B3" = "2.414"
    "C3" = "0.137"
    "D3" = "20.118"
    "E3" = "5.672"
    "F3" = "173.172"
    "G3" = "37.925"
    "H3" = "0.039"
    "I3" = "0.009"
    "J3" = "0.699"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}

Write-Host "Applied power-law model results"
